# Re-sort the AUC ranking rows (Dataset rows 4-8) alphabetically
# (ordinal / case-sensitive sort: uppercase before lowercase)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order (row -> data), matching a plain ordinal sort of dataset names
$rows = @(
    @{ Name = "CYP1A2"; Values = @(5, 6, 7, 8, 1, 4, 2, 2) },
    @{ Name = "Cancer"; Values = @(7, 5, 5, 8, 3, 3, 1, 2) },
    @{ Name = "HIV";    Values = @(4, 7, 8, 6, 1, 4, 2, 3) },
    @{ Name = "Liver";  Values = @(8, 5, 6, 7, 1, 3, 2, 3) },
    @{ Name = "hERG";   Values = @(5, 7, 8, 5, 1, 4, 2, 3) }
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i].Name
    $vals = $rows[$i].Values
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $vals[$c]
    }
}
